$wb = $excel.ActiveWorkbook

# OFF sheet - Week 15 "Road" (R) row updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 159
$wsOff.Range("C3").Value = 108
$wsOff.Range("D3").Value = 44
$wsOff.Range("E3").Value = 26

# DEF sheet - Week 15 "Road" (R) row updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 141
$wsDef.Range("C3").Value = 112
$wsDef.Range("D3").Value = 36
$wsDef.Range("E3").Value = 17
